$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.953.41"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "3.267.55"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.89"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.03"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "3.834.78"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.49"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "68.860.84"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").Value = "3.166.47"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.24"
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.91"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.192"
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.71"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.99"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.94"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.34"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.58"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("E43").Value = "  -5.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0690"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "345.79"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").Value = "2.603.06"
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.74"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.31"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.67"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  -0.65%  "
